$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 10
$ws.Range("B3").Value = "1E"
$ws.Range("A4").Value = 12
$ws.Range("B4").Value = "1E"

$ws.Range("A3:A4").HorizontalAlignment = -4131

$ws.Range("A4").Select()
